# Update "User Stories" worksheet:
#  - Row 3 (G3): status changes from text "In-progress" to a numeric time value
#    0.75 (18:00), matching the existing percentage/time style already used in G4.
#  - Row 6 (G6): status changes from "In-progress" to "Not yet started".
#  - New row 7 is added for a new user story about viewing a student's
#    time-in history (the "timeout" / time tracking feature).
#  - The active selection is moved to G7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: status is now a time value (75% of a day = 18:00) instead of text.
$ws.Range("G3").Value = 0.75

# Row 6: status text changes.
$ws.Range("G6").Value = "Not yet started"

# New row 7: additional user story.
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "user"
$ws.Range("C7").Value = "View a student's time in history"
$ws.Range("D7").Value = "I can see the activity of a student"
$ws.Range("F7").Value = "Normal"
$ws.Range("G7").Value = "Not yet started"

# Move the selection to the newly added cell, as in the final workbook.
$ws.Range("G7").Select() | Out-Null

# The real workbook also ends up scrolled one pane to the right (column D
# becomes the leftmost visible column) after navigating to G7. Attempt to
# reproduce that scroll position too; harmless if unsupported.
try {
    $win = $excel.ActiveWindow
    $win.SmallScroll($null, $null, 3, $null) | Out-Null
} catch {
}
